# Actualizacion automatica 2025-11-03 08:30:05
# Monthly rollover: clears the outgoing month's figures/"de 51" counters on
# "VENTAS POR GRUPO" and shifts "VENTA MENSUAL" one month to the right
# (julio/agosto/septiembre/octubre -> agosto/septiembre/octubre/noviembre).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("M9").Value = 0
$ws1.Range("M23").Value = 0
$ws1.Range("M32").Value = 0
$ws1.Range("R32").Value = 0
$ws1.Range("L39").Value = 0
$ws1.Range("D43").Value = 0
$ws1.Range("C47").Value = 0
$ws1.Range("C53").Value = "0 de 51"
$ws1.Range("D53").Value = "0 de 51"
$ws1.Range("L53").Value = "0 de 51"
$ws1.Range("M53").Value = "0 de 51"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C1").Value = "agosto"
$ws2.Range("D1").Value = "septiembre"
$ws2.Range("E1").Value = "octubre"
$ws2.Range("F1").Value = "noviembre"
$ws2.Range("C3").Value = 0
$ws2.Range("C4").Value = 0
$ws2.Range("C5").Value = 95.56
$ws2.Range("D5").Value = 0
$ws2.Range("C6").Value = 142.56
$ws2.Range("D6").Value = 0
$ws2.Range("D9").Value = 3864.45
$ws2.Range("E9").Value = -22.29
$ws2.Range("F9").Value = 0
$ws2.Range("D10").Value = 236.29
$ws2.Range("E10").Value = 0
$ws2.Range("D12").Value = 448.77
$ws2.Range("E12").Value = 0
$ws2.Range("D13").Value = 236.29
$ws2.Range("E13").Value = 0
$ws2.Range("C14").Value = 367.8
$ws2.Range("D14").Value = 0
$ws2.Range("D16").Value = 508.48
$ws2.Range("E16").Value = 0
$ws2.Range("C20").Value = 0
$ws2.Range("C21").Value = 354.43
$ws2.Range("D21").Value = 0
$ws2.Range("E23").Value = 56.02
$ws2.Range("F23").Value = 0
$ws2.Range("C27").Value = 0
$ws2.Range("D28").Value = 430.11
$ws2.Range("E28").Value = 0
$ws2.Range("D31").Value = 551.71
$ws2.Range("E31").Value = 0
$ws2.Range("E32").Value = -166.48
$ws2.Range("F32").Value = 0
$ws2.Range("C34").Value = 0
$ws2.Range("C35").Value = 0
$ws2.Range("C36").Value = 238.35
$ws2.Range("D36").Value = 0
$ws2.Range("C37").Value = 367.8
$ws2.Range("D37").Value = 0
$ws2.Range("D39").Value = 115.52
$ws2.Range("E39").Value = 179.12
$ws2.Range("F39").Value = 0
$ws2.Range("C41").Value = 0
$ws2.Range("D42").Value = 44.79
$ws2.Range("E42").Value = 0
$ws2.Range("E43").Value = 194.16
$ws2.Range("F43").Value = 0
$ws2.Range("C45").Value = 0
$ws2.Range("D46").Value = 5858.53
$ws2.Range("E46").Value = 0
$ws2.Range("E47").Value = 334.37
$ws2.Range("F47").Value = 0
$ws2.Range("D48").Value = -10.44
$ws2.Range("E48").Value = 0
$ws2.Range("C49").Value = 0
$ws2.Range("D49").Value = 522.8200000000001
$ws2.Range("E49").Value = 0
$ws2.Range("D50").Value = 165.83
$ws2.Range("E50").Value = 0
$ws2.Range("D52").Value = 438.86
$ws2.Range("E52").Value = 0
$ws2.Range("C53").Value = 1566.5
$ws2.Range("D53").Value = 13412.01
$ws2.Range("E53").Value = 574.9
$ws2.Range("F53").Value = 0

# --- Sheet 2 column widths ---
$ws2.Columns.Item(4).ColumnWidth = 15.166666666666666
$ws2.Columns.Item(5).ColumnWidth = 12.166666666666666
$ws2.Columns.Item(6).ColumnWidth = 14.166666666666666
